$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.808.64"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.933.45"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.05%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.47"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.32%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.19"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +9.60%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("D9").Value = "2.924.29"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("E13").Value = "  +3.95%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.77"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +5.02%  "
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "3.417.13"
$ws.Range("E16").Value = "  +1.43%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.89"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +9.25%  "
$ws.Range("D18").Value = "2.926.48"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "57.770.33"
$ws.Range("E19").Value = "  -0.32%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.55"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.50%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("E22").Value = "  +7.26%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.45"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +7.87%  "
$ws.Range("E24").Value = "  +3.24%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.28"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +6.34%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +5.68%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.38"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +3.70%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.94"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +0.78%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0979"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +5.73%  "
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("E35").Value = "  +6.06%  "
$ws.Range("E36").Value = "  +4.93%  "
$ws.Range("D37").Value = "0.0₃0695"
$ws.Range("E37").Value = "  +12.38%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.27"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +5.43%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.61"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +12.48%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.80%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "375.18"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +8.13%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0345"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "2.694.20"
$ws.Range("E44").Value = "  +3.46%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.08"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("E47").Value = "  +4.14%  "
$ws.Range("E48").Value = "  +2.43%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.95"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.98%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.90"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +3.34%  "
